# Added a README file with the project description
# (applies to the "To Do" sheet: mark the owner-edrpous item done, add a
#  reference link + two follow-up tasks)

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("To Do")

# --- Column B: add the reference link for row 1 ----------------------------
$ws2.Hyperlinks.Add($ws2.Range("B1"), "https://towardsdatascience.com/how-to-deploy-a-streamlit-app-using-an-amazon-free-ec2-instance-416a41f69dc3") | Out-Null
$ws2.Range("B2").Style = "Hyperlink"

# --- Column A: refresh the To Do list -------------------------------------
# Row 1 (Investigate how to deploy streamlit prototype to AWS) stays as-is.
# Row 2 used to be "Think what to do with owner edrpous" -> replaced with a
# brand new task. Row 3 is a new task appended below.
$ws2.Range("A2").Value = "How to update model with the new data"
$ws2.Range("A3").Value = "Code cleanup"

# --- Column widths -----------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 50.45
$ws2.Columns.Item(2).ColumnWidth = 95.6

# --- Selection ------------------------------------------------------------
$ws2.Range("A9").Select() | Out-Null
